$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.212.73"
$ws.Range("E2").Value = "  +2.81%  "
$ws.Range("D3").Value = "3.065.05"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'559.68"
$ws.Range("E5").Value = "  +2.58%  "
$ws.Range("D6").Value = "'143.71"
$ws.Range("E6").Value = "  +3.20%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.062.92"
$ws.Range("E8").Value = "  +2.13%  "
$ws.Range("E9").Value = "  +4.69%  "
$ws.Range("D10").Value = "'0.157"
$ws.Range("E10").Value = "  +5.57%  "
$ws.Range("D11").Value = "'6.12"
$ws.Range("E11").Value = "  -11.00%  "
$ws.Range("D12").Value = "'0.482"
$ws.Range("E12").Value = "  +7.76%  "
$ws.Range("E13").Value = "  +5.21%  "
$ws.Range("D14").Value = "'35.74"
$ws.Range("E14").Value = "  +4.85%  "
$ws.Range("D15").Value = "3.562.19"
$ws.Range("E15").Value = "  +2.37%  "
$ws.Range("D16").Value = "64.180.97"
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("D17").Value = "3.066.24"
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("D19").Value = "'6.78"
$ws.Range("E19").Value = "  +2.98%  "
$ws.Range("D20").Value = "'480.56"
$ws.Range("E20").Value = "  +1.86%  "
$ws.Range("D21").Value = "'14.02"
$ws.Range("E21").Value = "  +4.19%  "
$ws.Range("D22").Value = "'0.683"
$ws.Range("E22").Value = "  +4.37%  "
$ws.Range("D23").Value = "'14.52"
$ws.Range("E23").Value = "  +14.97%  "
$ws.Range("D24").Value = "'7.64"
$ws.Range("E24").Value = "  +6.12%  "
$ws.Range("D25").Value = "'82.24"
$ws.Range("E25").Value = "  +3.31%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").Value = "'2.81"
$ws.Range("E27").Value = "  +2.75%  "
$ws.Range("D28").Value = "'8.03"
$ws.Range("E28").Value = "  +4.95%  "
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").Value = "'26.43"
$ws.Range("E31").Value = "  +3.55%  "
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").Value = "'2.45"
$ws.Range("E33").Value = "  +4.08%  "
$ws.Range("D34").Value = "'5.71"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("D35").Value = "'6.25"
$ws.Range("E35").Value = "  +6.88%  "
$ws.Range("D36").Value = "'54.92"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").Value = "'0.0410"
$ws.Range("E37").Value = "  +4.41%  "
$ws.Range("D38").Value = "'446.79"
$ws.Range("E38").Value = "  -1.43%  "
$ws.Range("D39").Value = "'0.0814"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "'2.88"
$ws.Range("E40").Value = "  +11.52%  "
$ws.Range("D41").Value = "2.993.84"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").Value = "'8.28"
$ws.Range("E42").Value = "  +2.39%  "
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("D44").Value = "'28.06"
$ws.Range("E44").Value = "  +3.89%  "
$ws.Range("D45").Value = "'0.263"
$ws.Range("E45").Value = "  +5.63%  "
$ws.Range("D46").Value = "'2.17"
$ws.Range("E46").Value = "  +8.16%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "'0.114"
$ws.Range("E48").Value = "  +4.13%  "
$ws.Range("D49").Value = "0.0₃0521"
$ws.Range("E49").Value = "  +5.21%  "
$ws.Range("D50").Value = "'118.96"
$ws.Range("E50").Value = "  +3.07%  "
$ws.Range("E51").Value = "  +3.34%  "
